$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2158.1843
$ws.Range("I98").Value = 2259.0857
$ws.Range("K98").Value = 2259.0857
$ws.Range("M98").Value = -761.0857000000001
$ws.Range("H112").Value = 2215.5334
$ws.Range("J112").Value = 2215.5334
$ws.Range("L112").Value = 6646.600199999999
$ws.Range("N112").Value = -8862.600199999999
$ws.Range("H122").Value = 2158.1843
$ws.Range("I122").Value = 2259.0857
$ws.Range("K122").Value = 6777.257100000001
$ws.Range("M122").Value = -4327.257100000001
$ws.Range("H132").Value = 9208.667
$ws.Range("I132").Value = 9476.667
$ws.Range("J132").Value = 7332.6665
$ws.Range("K132").Value = 28430.001
$ws.Range("L132").Value = 21997.9995
$ws.Range("M132").Value = -25900.001
$ws.Range("N132").Value = -27057.9995
$ws.Range("H137").Value = 3127798.5
$ws.Range("I137").Value = 12501769
$ws.Range("J137").Value = 3141.5
$ws.Range("K137").Value = 37505307
$ws.Range("L137").Value = 9424.5
$ws.Range("M137").Value = -37502757
$ws.Range("N137").Value = -14524.5
$ws.Range("H138").Value = 3522.5195
$ws.Range("I138").Value = 3282.08
$ws.Range("J138").Value = 3638.1155
$ws.Range("K138").Value = 9846.24
$ws.Range("L138").Value = 10914.3465
$ws.Range("M138").Value = -4706.24
$ws.Range("N138").Value = -21194.3465
$ws.Range("H141").Value = 4032.7
$ws.Range("I141").Value = 1980.0385
$ws.Range("J141").Value = 17375
$ws.Range("K141").Value = 5940.1155
$ws.Range("L141").Value = 52125
$ws.Range("M141").Value = -760.1154999999999
$ws.Range("N141").Value = -62485

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2841.8696
$ws.Range("I32").Value = 2563.9062
$ws.Range("J32").Value = 6399.8
$ws.Range("K32").Value = 2563.9062
$ws.Range("L32").Value = 6399.8
$ws.Range("M32").Value = -2276.9062
$ws.Range("N32").Value = -6973.8
$ws.Range("H61").Value = 2384.0625
$ws.Range("I61").Value = 2016.8182
$ws.Range("J61").Value = 3192
$ws.Range("K61").Value = 2016.8182
$ws.Range("L61").Value = 3192
$ws.Range("M61").Value = -1804.8182
$ws.Range("N61").Value = -3616
$ws.Range("H110").Value = 4028.4119
$ws.Range("I110").Value = 2641.3333
$ws.Range("J110").Value = 7357.4
$ws.Range("K110").Value = 2641.3333
$ws.Range("L110").Value = 7357.4
$ws.Range("M110").Value = -596.3332999999998
$ws.Range("N110").Value = -11447.4
$ws.Range("H132").Value = 4446656
$ws.Range("I132").Value = 1709.2759
$ws.Range("J132").Value = 19611768
$ws.Range("K132").Value = 5127.8277
$ws.Range("L132").Value = 58835304
$ws.Range("M132").Value = -2597.8277
$ws.Range("N132").Value = -58840364
$ws.Range("H136").Value = 2384.0625
$ws.Range("I136").Value = 2016.8182
$ws.Range("J136").Value = 3192
$ws.Range("K136").Value = 6050.4546
$ws.Range("L136").Value = 9576
$ws.Range("M136").Value = -3500.4546
$ws.Range("N136").Value = -14676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 64517130
$ws.Range("I94").Value = 86957060
$ws.Range("J94").Value = 2327
$ws.Range("K94").Value = 86957060
$ws.Range("L94").Value = 2327
$ws.Range("M94").Value = -86956609
$ws.Range("N94").Value = -3229
$ws.Range("H134").Value = 2441.39
$ws.Range("I134").Value = 2185.5833
$ws.Range("J134").Value = 3557.6365
$ws.Range("K134").Value = 6556.749899999999
$ws.Range("L134").Value = 10672.9095
$ws.Range("M134").Value = -4021.749899999999
$ws.Range("N134").Value = -15742.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4608.96
$ws.Range("I31").Value = 3259.4285
$ws.Range("J31").Value = 5133.778
$ws.Range("K31").Value = 3259.4285
$ws.Range("L31").Value = 5133.778
$ws.Range("M31").Value = -2964.4285
$ws.Range("N31").Value = -5723.778
$ws.Range("H34").Value = 4608.96
$ws.Range("I34").Value = 3259.4285
$ws.Range("J34").Value = 5133.778
$ws.Range("K34").Value = 3259.4285
$ws.Range("L34").Value = 5133.778
$ws.Range("M34").Value = -3057.4285
$ws.Range("N34").Value = -5537.778
$ws.Range("H134").Value = 2925.6943
$ws.Range("I134").Value = 2953.1333
$ws.Range("J134").Value = 2788.5
$ws.Range("K134").Value = 8859.3999
$ws.Range("L134").Value = 8365.5
$ws.Range("M134").Value = -6324.3999
$ws.Range("N134").Value = -13435.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6754.273
$ws.Range("I131").Value = 14890.5
$ws.Range("J131").Value = 2105
$ws.Range("K131").Value = 44671.5
$ws.Range("L131").Value = 6315
$ws.Range("M131").Value = -39631.5
$ws.Range("N131").Value = -16395
$ws.Range("H140").Value = 9228.024
$ws.Range("I140").Value = 5807.448
$ws.Range("K140").Value = 17422.344
$ws.Range("M140").Value = -12242.344

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 114977.21
$ws.Range("I70").Value = 169797.33
$ws.Range("J70").Value = 20999.857
$ws.Range("K70").Value = 169797.33
$ws.Range("L70").Value = 20999.857
$ws.Range("M70").Value = -169527.33
$ws.Range("N70").Value = -21539.857
$ws.Range("H73").Value = 114977.21
$ws.Range("I73").Value = 169797.33
$ws.Range("J73").Value = 20999.857
$ws.Range("K73").Value = 169797.33
$ws.Range("L73").Value = 20999.857
$ws.Range("M73").Value = -168861.33
$ws.Range("N73").Value = -22871.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 496.75
$ws.Range("I16").Value = 496.75
$ws.Range("K16").Value = 496.75
$ws.Range("M16").Value = -326.75
$ws.Range("H82").Value = 1307.3954
$ws.Range("I82").Value = 1318.6389
$ws.Range("J82").Value = 1249.5714
$ws.Range("K82").Value = 1318.6389
$ws.Range("L82").Value = 1249.5714
$ws.Range("M82").Value = -957.6388999999999
$ws.Range("N82").Value = -1971.5714
$ws.Range("H85").Value = 1307.3954
$ws.Range("I85").Value = 1318.6389
$ws.Range("J85").Value = 1249.5714
$ws.Range("K85").Value = 1318.6389
$ws.Range("L85").Value = 1249.5714
$ws.Range("M85").Value = -70.63889999999992
$ws.Range("N85").Value = -3745.5714
$ws.Range("H100").Value = 1649
$ws.Range("I100").Value = 1399.5
$ws.Range("J100").Value = 1898.5
$ws.Range("K100").Value = 1399.5
$ws.Range("L100").Value = 1898.5
$ws.Range("M100").Value = -858.5
$ws.Range("N100").Value = -2980.5
$ws.Range("H132").Value = 3690.2903
$ws.Range("I132").Value = 3278.889
$ws.Range("J132").Value = 4259.923
$ws.Range("K132").Value = 9836.667000000001
$ws.Range("L132").Value = 12779.769
$ws.Range("M132").Value = -7306.667000000001
$ws.Range("N132").Value = -17839.769
$ws.Range("H136").Value = 2764.6858
$ws.Range("I136").Value = 2621.4194
$ws.Range("K136").Value = 7864.2582
$ws.Range("M136").Value = -5314.2582

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 43478976
$ws.Range("I100").Value = 596.2353
$ws.Range("K100").Value = 1192.4706
$ws.Range("M100").Value = -651.4706000000001
$ws.Range("H132").Value = 12349876
$ws.Range("I132").Value = 17547968
$ws.Range("J132").Value = 4405.5
$ws.Range("K132").Value = 52643904
$ws.Range("L132").Value = 13216.5
$ws.Range("M132").Value = -52641374
$ws.Range("N132").Value = -18276.5
$ws.Range("H136").Value = 5381.378
$ws.Range("I136").Value = 5372.1284
$ws.Range("K136").Value = 16116.3852
$ws.Range("M136").Value = -13566.3852
